# "little changes to final ppt v3"
# 1) Slide 40: "normally distributed " -> "uniformly distributed " (bold run)
# 2) Slide 41: "...were normally distributed..." -> "...were uniformly distributed..."
# 3) Slide 42: merge " FSM coverage has a low " + "resul" + "." runs into
#    a single run reading " FSM coverage has a low result."

$p = $ppt.ActivePresentation

# --- Change 1: slide 40 ---
$s40 = $p.Slides.Item(40)
$sh40 = $s40.Shapes.Item(2)
$tr40 = $sh40.TextFrame.TextRange
$run40 = $tr40.Characters(59, 21)
$run40.Text = "uniformly distributed "

# --- Change 2: slide 41 ---
$s41 = $p.Slides.Item(41)
$sh41 = $s41.Shapes.Item(2)
$tr41 = $sh41.TextFrame.TextRange
$run41 = $tr41.Characters(16, 106)
$run41.Text = "The values of its coordinates were uniformly distributed between all the possible values in all the tests."

# --- Change 3: slide 42 ---
$s42 = $p.Slides.Item(42)
$sh42 = $s42.Shapes.Item(2)
$tr42 = $sh42.TextFrame.TextRange
$run42 = $tr42.Characters(12, 30)
$run42.Text = " FSM coverage has a low result."
